# Scheduled runner update: refresh Leve profit calculations (current market
# prices, leve rewards, and resulting NQ/HQ profit margins) across the
# crafting-job worksheets, per the latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1120.4584
$ws.Range("I18").Value = 1156.1305
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 1156.1305
$ws.Range("L18").Value = 300
$ws.Range("M18").Value = -872.1305
$ws.Range("N18").Value = -868
# Row 39
$ws.Range("H39").Value = 21740856
$ws.Range("I39").Value = 55.333332
$ws.Range("J39").Value = 62504856
$ws.Range("K39").Value = 165.999996
$ws.Range("L39").Value = 187514568
$ws.Range("M39").Value = 130.000004
$ws.Range("N39").Value = -187515160
# Row 137
$ws.Range("H137").Value = 5245676
$ws.Range("I137").Value = 6692948
$ws.Range("J137").Value = 4184343.2
$ws.Range("K137").Value = 20078844
$ws.Range("L137").Value = 12553029.6
$ws.Range("M137").Value = -20076294
$ws.Range("N137").Value = -12558129.6
# Row 141
$ws.Range("H141").Value = 7439.7
$ws.Range("I141").Value = 4230
$ws.Range("J141").Value = 8815.286
$ws.Range("K141").Value = 12690
$ws.Range("L141").Value = 26445.858
$ws.Range("M141").Value = -7510
$ws.Range("N141").Value = -36805.858

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1512892.9
$ws.Range("I32").Value = 1651337
$ws.Range("K32").Value = 1651337
$ws.Range("M32").Value = -1651050
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 131
$ws.Range("H131").Value = 58986.668
$ws.Range("J131").Value = 58986.668
$ws.Range("L131").Value = 58986.668
$ws.Range("N131").Value = -69066.66800000001

$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Range("H60").Value = 23940
$ws.Range("J60").Value = 23940
$ws.Range("L60").Value = 23940
$ws.Range("N60").Value = -25138
# Row 86
$ws.Range("H86").Value = 1876.6566
$ws.Range("I86").Value = 1911.978
$ws.Range("J86").Value = 1474.875
$ws.Range("K86").Value = 1911.978
$ws.Range("L86").Value = 1474.875
$ws.Range("M86").Value = -788.9780000000001
$ws.Range("N86").Value = -3720.875
# Row 89
$ws.Range("H89").Value = 1876.6566
$ws.Range("I89").Value = 1911.978
$ws.Range("J89").Value = 1474.875
$ws.Range("K89").Value = 9559.889999999999
$ws.Range("L89").Value = 7374.375
$ws.Range("M89").Value = -3943.889999999999
$ws.Range("N89").Value = -18606.375

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1538706.9
$ws.Range("I22").Value = 1818408.1
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 1818408.1
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -1818058.1
$ws.Range("N22").Value = -1050
# Row 62
$ws.Range("H62").Value = 2442.8572
$ws.Range("I62").Value = 2442.8572
$ws.Range("K62").Value = 2442.8572
$ws.Range("M62").Value = -1818.8572
# Row 65
$ws.Range("H65").Value = 2442.8572
$ws.Range("I65").Value = 2442.8572
$ws.Range("K65").Value = 12214.286
$ws.Range("M65").Value = -9094.286
# Row 131
$ws.Range("H131").Value = 15320
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 15320
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 15320
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -25400

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 102.1875
$ws.Range("I26").Value = 108.28571
$ws.Range("J26").Value = 59.5
$ws.Range("K26").Value = 324.85713
$ws.Range("L26").Value = 178.5
$ws.Range("M26").Value = -36.85712999999998
$ws.Range("N26").Value = -754.5
# Row 68
$ws.Range("H68").Value = 2467.4385
$ws.Range("I68").Value = 671.3214
$ws.Range("J68").Value = 3585.0222
$ws.Range("K68").Value = 2013.9642
$ws.Range("L68").Value = 10755.0666
$ws.Range("M68").Value = -1202.9642
$ws.Range("N68").Value = -12377.0666
# Row 69
$ws.Range("H69").Value = 1213.3684
$ws.Range("I69").Value = 928.4286
$ws.Range("J69").Value = 1379.5834
$ws.Range("K69").Value = 2785.2858
$ws.Range("L69").Value = 4138.7502
$ws.Range("M69").Value = -1974.2858
$ws.Range("N69").Value = -5760.7502
# Row 71
$ws.Range("H71").Value = 2467.4385
$ws.Range("I71").Value = 671.3214
$ws.Range("J71").Value = 3585.0222
$ws.Range("K71").Value = 6041.8926
$ws.Range("L71").Value = 32265.1998
$ws.Range("M71").Value = -1985.8926
$ws.Range("N71").Value = -40377.1998
# Row 72
$ws.Range("H72").Value = 1213.3684
$ws.Range("I72").Value = 928.4286
$ws.Range("J72").Value = 1379.5834
$ws.Range("K72").Value = 8355.857399999999
$ws.Range("L72").Value = 12416.2506
$ws.Range("M72").Value = -4299.857399999999
$ws.Range("N72").Value = -20528.2506
# Row 76
$ws.Range("H76").Value = 5582.4
$ws.Range("I76").Value = 2970.6667
$ws.Range("J76").Value = 9500
$ws.Range("K76").Value = 8912.000100000001
$ws.Range("L76").Value = 28500
$ws.Range("M76").Value = -8529.000100000001
$ws.Range("N76").Value = -29266
# Row 79
$ws.Range("H79").Value = 5582.4
$ws.Range("I79").Value = 2970.6667
$ws.Range("J79").Value = 9500
$ws.Range("K79").Value = 8912.000100000001
$ws.Range("L79").Value = 28500
$ws.Range("M79").Value = -7586.000100000001
$ws.Range("N79").Value = -31152

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3889920.2
$ws.Range("I122").Value = 27517.2
$ws.Range("J122").Value = 7248531.5
$ws.Range("K122").Value = 82551.60000000001
$ws.Range("L122").Value = 21745594.5
$ws.Range("M122").Value = -80101.60000000001
$ws.Range("N122").Value = -21750494.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1368.4445
$ws.Range("I46").Value = 2666.5
$ws.Range("J46").Value = 330
$ws.Range("K46").Value = 2666.5
$ws.Range("L46").Value = 330
$ws.Range("M46").Value = -2478.5
$ws.Range("N46").Value = -706
# Row 68
$ws.Range("H68").Value = 2154.5
$ws.Range("I68").Value = 1964.0625
$ws.Range("J68").Value = 2535.375
$ws.Range("K68").Value = 1964.0625
$ws.Range("L68").Value = 2535.375
$ws.Range("M68").Value = -1215.0625
$ws.Range("N68").Value = -4033.375
# Row 71
$ws.Range("H71").Value = 2154.5
$ws.Range("I71").Value = 1964.0625
$ws.Range("J71").Value = 2535.375
$ws.Range("K71").Value = 9820.3125
$ws.Range("L71").Value = 12676.875
$ws.Range("M71").Value = -6076.3125
$ws.Range("N71").Value = -20164.875
# Row 104
$ws.Range("H104").Value = 18266
$ws.Range("J104").Value = 18266
$ws.Range("L104").Value = 18266
$ws.Range("N104").Value = -25254

$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 8663.166999999999
$ws.Range("J101").Value = 8663.166999999999
$ws.Range("L101").Value = 8663.166999999999
$ws.Range("N101").Value = -15153.167
